$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing validation values (B2:B10)
$ws.Range("B2").Value = 0.2583083057973862
$ws.Range("B3").Value = 0.4551466506160821
$ws.Range("B4").Value = 0.08067549060998827
$ws.Range("B5").Value = 0.4721139708585126
$ws.Range("B6").Value = 0.2081995091274567
$ws.Range("B7").Value = 0.6412846493632511
$ws.Range("B8").Value = 0.1210795967132373
$ws.Range("B9").Value = 0.8745659912710601
$ws.Range("B10").Value = 0.1850681536555143

# Add new rows 11 and 12
$ws.Range("A11").Value = "6_1"
$ws.Range("B11").Value = 0.8564738596801792
$ws.Range("A12").Value = "6_2"
$ws.Range("B12").Value = 0.06405097250167674
